$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 7324.2
$ws.Cells.Item(86, 9).Value = 5075
$ws.Cells.Item(86, 10).Value = 13509.5
$ws.Cells.Item(86, 11).Value = 5075
$ws.Cells.Item(86, 12).Value = 13509.5
$ws.Cells.Item(86, 13).Value = -3952
$ws.Cells.Item(86, 14).Value = -15755.5
$ws.Cells.Item(89, 8).Value = 7324.2
$ws.Cells.Item(89, 9).Value = 5075
$ws.Cells.Item(89, 10).Value = 13509.5
$ws.Cells.Item(89, 11).Value = 25375
$ws.Cells.Item(89, 12).Value = 67547.5
$ws.Cells.Item(89, 13).Value = -19759
$ws.Cells.Item(89, 14).Value = -78779.5
$ws.Cells.Item(92, 8).Value = 679.9375
$ws.Cells.Item(92, 9).Value = 742.61536
$ws.Cells.Item(92, 10).Value = 408.33334
$ws.Cells.Item(92, 11).Value = 742.61536
$ws.Cells.Item(92, 12).Value = 408.33334
$ws.Cells.Item(92, 13).Value = 505.38464
$ws.Cells.Item(92, 14).Value = -2904.33334
$ws.Cells.Item(106, 8).Value = 2290.3
$ws.Cells.Item(106, 9).Value = 2378.111
$ws.Cells.Item(106, 11).Value = 2378.111
$ws.Cells.Item(106, 13).Value = -1747.111
$ws.Cells.Item(138, 8).Value = 6318.697
$ws.Cells.Item(138, 10).Value = 6980.8
$ws.Cells.Item(138, 12).Value = 20942.4
$ws.Cells.Item(138, 14).Value = -31222.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21511.576
$ws.Cells.Item(32, 9).Value = 3796.6233
$ws.Cells.Item(32, 10).Value = 97907.31
$ws.Cells.Item(32, 11).Value = 3796.6233
$ws.Cells.Item(32, 12).Value = 97907.31
$ws.Cells.Item(32, 13).Value = -3509.6233
$ws.Cells.Item(32, 14).Value = -98481.31
$ws.Cells.Item(61, 8).Value = 1497.4656
$ws.Cells.Item(61, 9).Value = 672.26086
$ws.Cells.Item(61, 10).Value = 2039.7428
$ws.Cells.Item(61, 11).Value = 672.26086
$ws.Cells.Item(61, 12).Value = 2039.7428
$ws.Cells.Item(61, 13).Value = -460.26086
$ws.Cells.Item(61, 14).Value = -2463.7428
$ws.Cells.Item(122, 8).Value = 3048.0356
$ws.Cells.Item(122, 9).Value = 3406.1765
$ws.Cells.Item(122, 10).Value = 2494.5454
$ws.Cells.Item(122, 11).Value = 10218.5295
$ws.Cells.Item(122, 12).Value = 7483.6362
$ws.Cells.Item(122, 13).Value = -7768.529500000001
$ws.Cells.Item(122, 14).Value = -12383.6362
$ws.Cells.Item(136, 8).Value = 1497.4656
$ws.Cells.Item(136, 9).Value = 672.26086
$ws.Cells.Item(136, 10).Value = 2039.7428
$ws.Cells.Item(136, 11).Value = 2016.78258
$ws.Cells.Item(136, 12).Value = 6119.2284
$ws.Cells.Item(136, 13).Value = 533.2174199999999
$ws.Cells.Item(136, 14).Value = -11219.2284
$ws.Cells.Item(137, 8).Value = 39666.668
$ws.Cells.Item(137, 10).Value = 39666.668
$ws.Cells.Item(137, 12).Value = 39666.668
$ws.Cells.Item(137, 14).Value = -49866.668
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 13990.378
$ws.Cells.Item(31, 9).Value = 28852.25
$ws.Cells.Item(31, 10).Value = 2359.348
$ws.Cells.Item(31, 11).Value = 28852.25
$ws.Cells.Item(31, 12).Value = 2359.348
$ws.Cells.Item(31, 13).Value = -28557.25
$ws.Cells.Item(31, 14).Value = -2949.348
$ws.Cells.Item(34, 8).Value = 13990.378
$ws.Cells.Item(34, 9).Value = 28852.25
$ws.Cells.Item(34, 10).Value = 2359.348
$ws.Cells.Item(34, 11).Value = 28852.25
$ws.Cells.Item(34, 12).Value = 2359.348
$ws.Cells.Item(34, 13).Value = -28650.25
$ws.Cells.Item(34, 14).Value = -2763.348
$ws.Cells.Item(132, 8).Value = 2953.4285
$ws.Cells.Item(132, 10).Value = 2913.8572
$ws.Cells.Item(132, 12).Value = 8741.571599999999
$ws.Cells.Item(132, 14).Value = -13801.5716

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2115.0598
$ws.Cells.Item(68, 9).Value = 1403.3214
$ws.Cells.Item(68, 11).Value = 4209.9642
$ws.Cells.Item(68, 13).Value = -3398.9642
$ws.Cells.Item(71, 8).Value = 2115.0598
$ws.Cells.Item(71, 9).Value = 1403.3214
$ws.Cells.Item(71, 11).Value = 12629.8926
$ws.Cells.Item(71, 13).Value = -8573.892600000001
$ws.Cells.Item(113, 8).Value = 550.76746
$ws.Cells.Item(113, 10).Value = 558.8214
$ws.Cells.Item(113, 12).Value = 1676.4642
$ws.Cells.Item(113, 14).Value = -6016.4642

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2003.409
$ws.Cells.Item(102, 9).Value = 2109.3157
$ws.Cells.Item(102, 10).Value = 1332.6666
$ws.Cells.Item(102, 11).Value = 2109.3157
$ws.Cells.Item(102, 12).Value = 1332.6666
$ws.Cells.Item(102, 13).Value = -487.3157000000001
$ws.Cells.Item(102, 14).Value = -4576.6666
$ws.Cells.Item(126, 8).Value = 1695.0869
$ws.Cells.Item(126, 9).Value = 1563.9412
$ws.Cells.Item(126, 10).Value = 2066.6667
$ws.Cells.Item(126, 11).Value = 4691.8236
$ws.Cells.Item(126, 12).Value = 6200.000100000001
$ws.Cells.Item(126, 13).Value = -2221.8236
$ws.Cells.Item(126, 14).Value = -11140.0001
$ws.Cells.Item(135, 8).Value = 44090.77
$ws.Cells.Item(135, 10).Value = 44090.77
$ws.Cells.Item(135, 12).Value = 44090.77
$ws.Cells.Item(135, 14).Value = -54230.77
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(56, 8).Value = 33723.168
$ws.Cells.Item(56, 10).Value = 36490.273
$ws.Cells.Item(56, 12).Value = 36490.273
$ws.Cells.Item(56, 14).Value = -37918.273
$ws.Cells.Item(81, 8).Value = 250775.12
$ws.Cells.Item(81, 9).Value = 200940
$ws.Cells.Item(81, 10).Value = 333833.66
$ws.Cells.Item(81, 11).Value = 401880
$ws.Cells.Item(81, 12).Value = 667667.3199999999
$ws.Cells.Item(81, 13).Value = -400819
$ws.Cells.Item(81, 14).Value = -669789.3199999999
$ws.Cells.Item(84, 8).Value = 250775.12
$ws.Cells.Item(84, 9).Value = 200940
$ws.Cells.Item(84, 10).Value = 333833.66
$ws.Cells.Item(84, 11).Value = 2009400
$ws.Cells.Item(84, 12).Value = 3338336.6
$ws.Cells.Item(84, 13).Value = -2004096
$ws.Cells.Item(84, 14).Value = -3348944.6
$ws.Cells.Item(132, 8).Value = 3126.0303
$ws.Cells.Item(132, 9).Value = 3950.5789
$ws.Cells.Item(132, 10).Value = 2007
$ws.Cells.Item(132, 11).Value = 11851.7367
$ws.Cells.Item(132, 12).Value = 6021
$ws.Cells.Item(132, 13).Value = -9321.736699999999
$ws.Cells.Item(132, 14).Value = -11081
$ws.Cells.Item(136, 8).Value = 1323.6
$ws.Cells.Item(136, 9).Value = 789.7692
$ws.Cells.Item(136, 10).Value = 2315
$ws.Cells.Item(136, 11).Value = 2369.3076
$ws.Cells.Item(136, 12).Value = 6945
$ws.Cells.Item(136, 13).Value = 180.6923999999999
$ws.Cells.Item(136, 14).Value = -12045
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()
